# Update market health data on "Top Gainers" sheet
# Commit: Update data: 2025-10-30 13:32
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Top Gainers")

$rows = @(
    @(2, "ASALCBR", 13.7911, 14.0037, 28.4881),
    @(3, "VAIBHAVGBL", 11.1823, 19.9891, 27.1472),
    @(4, "SALASAR", 11.066, 16.383, 23.3371),
    @(5, "KELLTONTEC", 9.0909, 6.174, -8.5044),
    @(6, "SAGILITY", 8.787100000000001, 16.4562, 29.7234),
    @(7, "UNIPARTS", 8.491, 10.7844, 27.0959),
    @(8, "CHENNPETRO", 7.6201, 12.8245, 14.9908),
    @(9, "SHREEJISPG", 6.8483, 11.1523, 12.2547),
    @(10, "MARINE", 6.5114, 3.035, 15.6171),
    @(11, "POLICYBZR", 6.0868, 8.4594, 7.4207),
    @(12, "MCLOUD", 5.9319, 4.9166, -23.2437),
    @(13, "BLSE", 5.6317, 4.5958, -1.6103),
    @(14, "INDIACEM", 5.6088, 5.8929, 7.6302),
    @(15, "PDSL", 5.594, 8.939500000000001, 15.1029),
    @(16, "VENKEYS", 5.4233, 6.0527, 3.746),
    @(17, "BHEL", 5.4036, 11.9988, 8.471399999999999),
    @(18, "NETWEB", 5.3718, 11.2426, 13.3342),
    @(19, "SOLEX", 4.9978, 6.5692, "N/A"),
    @(20, "V2RETAIL", 4.9975, 3.4793, 11.7421),
    @(21, "HIRECT", 4.9951, 12.5078, 10.7801),
    @(22, "SKMEGGPROD", 4.9909, 12.0154, 29.9408),
    @(23, "EUROPRATIK", 4.9438, 10.8518, 27.9093),
    @(24, "VIMTALABS", 4.88, 5.0595, -0.0708),
    @(25, "RAMASTEEL", 4.4177, 4.3129, 5.9063),
    @(26, "REDTAPE", 4.4083, 4.3069, -2.6745),
    @(27, "MTARTECH", 4.3689, 8.5245, 32.6456),
    @(28, "SKYGOLD", 4.3628, -0.2188, 38.5755),
    @(29, "BAJAJHCARE", 4.2697, 4.8177, -1.4718),
    @(30, "CREDITACC", 4.1604, 2.7838, 8.047599999999999),
    @(31, "ABREL", 4.0653, 12.4406, 11.9747),
    @(32, "DEEDEV", 4.0494, -2.8728, -3.6738),
    @(33, "ALICON", 4.0406, 10.2083, 15.6826),
    @(34, "IVALUE", 3.9881, 7.2664, 0),
    @(35, "SHANTIGOLD", 3.8985, 11.2113, 3.7917),
    @(36, "DBL", 3.8744, 4.9728, 6.0609),
    @(37, "BLUEDART", 3.8725, 22.6574, 19.7002),
    @(38, "PFOCUS", 3.8493, 1.1205, 2.5861),
    @(39, "SHRINGARMS", 3.8237, 5.0379, 25.3279),
    @(40, "CANBK", 3.6422, 6.1734, 7.8813),
    @(41, "RSYSTEMS", 3.6318, 4.795, 7.1482),
    @(42, "AHLUCONT", 3.6218, 2.4227, -4.8432),
    @(43, "CENTRUM", 3.6061, 2.2123, 1.5444),
    @(44, "SAPPHIRE", 3.5684, 5.3946, 2.74),
    @(45, "BGRENERGY", 3.5346, -5.8011, 75.3689),
    @(46, "RELTD", 3.4551, 10.2584, -1.2921),
    @(47, "BLISSGVS", 3.4528, 2.8037, 3.1746),
    @(48, "VSTIND", 3.2682, 3.7699, 3.3082),
    @(49, "INOXGREEN", 3.1231, 11.1631, 34.6328),
    @(50, "INDORAMA", 3.0612, 5.5381, 16.9755),
    @(51, "FIVESTAR", 3.0576, 16.1173, 16.2039),
    @(52, "NEULANDLAB", 3.0443, -1.324, 8.6957),
    @(53, "GRMOVER", 3.0307, 3.2051, 18.8937),
    @(54, "SUNDROP", 2.9786, 2.8713, 0.9714),
    @(55, "POWERINDIA", 2.9573, 7.3478, -0.0111),
    @(56, "ASHOKA", 2.8201, 4.3403, 6.9665),
    @(57, "MRPL", 2.8165, 12.8002, 23.4355),
    @(58, "GMMPFAUDLR", 2.8152, 7.271, 19.5905),
    @(59, "PSPPROJECT", 2.7993, 16.8631, 23.2881),
    @(60, "OIL", 2.7705, 3.0155, 4.4345),
    @(61, "VOLTAMP", 2.6913, 2.5703, 2.3081),
    @(62, "BPCL", 2.6716, 8.1555, 5.226),
    @(63, "REFEX", 2.6657, 0.0963, 2.0755),
    @(64, "MFSL", 2.6337, 2.6878, -1.0851),
    @(65, "JKLAKSHMI", 2.5851, 4.5741, 1.592),
    @(66, "JKTYRE", 2.5816, 5.6196, 21.6553),
    @(67, "BLS", 2.5579, -0.4713, -1.722),
    @(68, "AARTIDRUGS", 2.5306, 2.5306, 4.1007),
    @(69, "CIFL", 2.5007, 2.0556, 1.9965),
    @(70, "CAMS", 2.4899, 1.8375, 5.1322),
    @(71, "WEBELSOLAR", 2.4636, 2.86, -1.0039),
    @(72, "SDBL", 2.46, 1.0066, 6.591),
    @(73, "IIFL", 2.4474, 9.3916, 18.5661),
    @(74, "WESTLIFE", 2.4091, 2.3656, -12.5218),
    @(75, "KMEW", 2.4057, 4.5145, -1.0415),
    @(76, "NBCC", 2.352, 5.5869, 10.1326)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

Write-Host "Applied market health data update to $($rows.Count) rows"
